$d = $word.ActiveDocument

# Locate the target paragraph (the one beginning with "Vous allez participer...").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Vous allez participer*") {
        $target = $p
        break
    }
}

$r = $target.Range

# Apply the wording change within the paragraph text: the constellation name
# changes from "constellation Persée" to "Constellation d'Hercule".
$oldPhrase = "constellation Persée"
$newPhrase = "Constellation d'Hercule"
$fullText = $r.Text
if ($fullText.Length -gt 0 -and $fullText.Substring($fullText.Length - 1) -eq [char]13) {
    $fullText = $fullText.Substring(0, $fullText.Length - 1)
}
$newText = $fullText.Replace($oldPhrase, $newPhrase)

# Restrict the range to the paragraph's text only (exclude the trailing
# paragraph mark) so the paragraph break/following empty paragraph survive.
$textRange = $d.Range($r.Start, $r.End - 1)

# Replace the whole paragraph content with a single freshly-typed run that
# carries no explicit run formatting (matching a full retype of the line).
$textRange.Delete()
$d.Range($r.Start, $r.Start).InsertBefore($newText)
